$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet, positioned right before "总计".
#    We copy the "2021-Q4" sheet (same A1:H10 shape, same column
#    layout/types) so the new sheet starts out with the correct
#    sheetPr / column typing, then we overwrite every value.
#    NOTE: sheet variables in this runtime are resolved by their
#    *current* position, so after the copy/insert we must re-fetch
#    any sheet object whose index shifted (e.g. "总计" moves from
#    index 6 to index 7).
# ------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(5)
$beforeSheet = $wb.Worksheets.Item(6)
$templateSheet.Copy($beforeSheet)
$newSheet = $wb.Worksheets.Item(6)
$newSheet.Name = "2022-Q1"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$data = @(
    @(0, "166005", "中欧价值发现混合 -A", "43.52", "93.97", "4.14", "1.8017", 5),
    @(1, "001882", "中欧价值发现混合 -E", "43.52", "93.97", "4.14", "1.8017", 5),
    @(2, "001810", "中欧潜力价值灵活配置混合A", "28.67", "94.05", "4.32", "1.2385", 4),
    @(3, "004232", "中欧价值发现混合 -C", "10.98", "93.97", "4.14", "0.4546", 5),
    @(4, "166024", "中欧恒利三年定期开放混合", "4.48", "98.71", "4.78", "0.2141", 5),
    @(5, "005764", "中欧潜力价值灵活配置混合C", "3.43", "94.05", "4.32", "0.1482", 4),
    @(6, "001891", "中欧成长优选回报灵活配置混合E", "2.97", "94.42", "3.48", "0.1034", 5),
    @(7, "166020", "中欧成长优选回报灵活配置混合A", "2.97", "94.42", "3.48", "0.1034", 5),
    @(8, "000963", "兴业多策略灵活配置混合", "2.07", "75.34", "3.81", "0.0789", 4)
)

$row = 2
foreach ($d in $data) {
    $newSheet.Range("A$row").Value = $d[0]
    $newSheet.Range("B$row").Value = "'" + $d[1]
    $newSheet.Range("C$row").Value = $d[2]
    $newSheet.Range("D$row").Value = "'" + $d[3]
    $newSheet.Range("E$row").Value = "'" + $d[4]
    $newSheet.Range("F$row").Value = "'" + $d[5]
    $newSheet.Range("G$row").Value = "'" + $d[6]
    $newSheet.Range("H$row").Value = $d[7]
    $row++
}

# ------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new row for 2022-Q1 at the
#    top of the data (row 2), push the rest down, renumber index
#    column, and fill in the new row's values.
#    "总计" is now at index 7 (it shifted when "2022-Q1" was added),
#    so fetch it fresh rather than reusing an old reference.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(7)
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 5.94

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
